$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 457.7857
$ws.Range("J17").Value = 457.7857
$ws.Range("L17").Value = 1373.3571
$ws.Range("N17").Value = -1709.3571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3085.5715
$ws.Range("I43").Value = 3085.5715
$ws.Range("K43").Value = 3085.5715
$ws.Range("M43").Value = -3016.5715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3397
$ws.Range("I106").Value = 3579.6667
$ws.Range("J106").Value = 2575
$ws.Range("K106").Value = 3579.6667
$ws.Range("L106").Value = 2575
$ws.Range("M106").Value = -2948.6667
$ws.Range("N106").Value = -3837

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2858.7856
$ws.Range("I113").Value = 2670.3333
$ws.Range("J113").Value = 2910.182
$ws.Range("K113").Value = 2670.3333
$ws.Range("L113").Value = 2910.182
$ws.Range("M113").Value = 583.6667000000002
$ws.Range("N113").Value = -9418.182000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6764.05
$ws.Range("I132").Value = 7769.5293
$ws.Range("K132").Value = 23308.5879
$ws.Range("M132").Value = -20778.5879

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2719.6511
$ws.Range("I137").Value = 1413.75
$ws.Range("J137").Value = 5157.3335
$ws.Range("K137").Value = 4241.25
$ws.Range("L137").Value = 15472.0005
$ws.Range("M137").Value = -1691.25
$ws.Range("N137").Value = -20572.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4565.0835
$ws.Range("I138").Value = 2794.8
$ws.Range("J138").Value = 5245.9614
$ws.Range("K138").Value = 8384.400000000001
$ws.Range("L138").Value = 15737.8842
$ws.Range("M138").Value = -3244.400000000001
$ws.Range("N138").Value = -26017.8842

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 1872726.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2655.8
$ws.Range("I141").Value = 2655.8
$ws.Range("K141").Value = 7967.400000000001
$ws.Range("M141").Value = -2787.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1870.3914
$ws.Range("I2").Value = 1073
$ws.Range("K2").Value = 1073
$ws.Range("M2").Value = -960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1969.4694
$ws.Range("I32").Value = 2051.4092
$ws.Range("K32").Value = 2051.4092
$ws.Range("M32").Value = -1764.4092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5185.839
$ws.Range("I61").Value = 4339.524
$ws.Range("J61").Value = 6963.1
$ws.Range("K61").Value = 4339.524
$ws.Range("L61").Value = 6963.1
$ws.Range("M61").Value = -4127.524
$ws.Range("N61").Value = -7387.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 29414114
$ws.Range("I74").Value = 1892.375
$ws.Range("K74").Value = 1892.375
$ws.Range("M74").Value = -1018.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 29414114
$ws.Range("I77").Value = 1892.375
$ws.Range("K77").Value = 9461.875
$ws.Range("M77").Value = -5093.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2473.8667
$ws.Range("I102").Value = 2220.4
$ws.Range("J102").Value = 2980.8
$ws.Range("K102").Value = 2220.4
$ws.Range("L102").Value = 2980.8
$ws.Range("M102").Value = -598.4000000000001
$ws.Range("N102").Value = -6224.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1548.375
$ws.Range("I110").Value = 1507
$ws.Range("K110").Value = 1507
$ws.Range("M110").Value = 538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1870.3914
$ws.Range("I116").Value = 1073
$ws.Range("K116").Value = 1073
$ws.Range("M116").Value = 1221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5185.839
$ws.Range("I136").Value = 4339.524
$ws.Range("J136").Value = 6963.1
$ws.Range("K136").Value = 13018.572
$ws.Range("L136").Value = 20889.3
$ws.Range("M136").Value = -10468.572
$ws.Range("N136").Value = -25989.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1870.3914
$ws.Range("I3").Value = 1073
$ws.Range("K3").Value = 1073
$ws.Range("M3").Value = -959

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3270.875
$ws.Range("I99").Value = 2444.5
$ws.Range("J99").Value = 5750
$ws.Range("K99").Value = 2444.5
$ws.Range("L99").Value = 5750
$ws.Range("M99").Value = -946.5
$ws.Range("N99").Value = -8746

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15295506
$ws.Range("I105").Value = 1251340.6
$ws.Range("J105").Value = 27779208
$ws.Range("K105").Value = 1251340.6
$ws.Range("L105").Value = 27779208
$ws.Range("M105").Value = -1249593.6
$ws.Range("N105").Value = -27782702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3498127.8
$ws.Range("I107").Value = 4274903
$ws.Range("K107").Value = 4274903
$ws.Range("M107").Value = -4272983

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2947.5605
$ws.Range("I31").Value = 1320.5
$ws.Range("K31").Value = 1320.5
$ws.Range("M31").Value = -1025.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2947.5605
$ws.Range("I34").Value = 1320.5
$ws.Range("K34").Value = 1320.5
$ws.Range("M34").Value = -1118.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2691.7896
$ws.Range("I58").Value = 1844.3
$ws.Range("J58").Value = 3633.4443
$ws.Range("K58").Value = 1844.3
$ws.Range("L58").Value = 3633.4443
$ws.Range("M58").Value = -1641.3
$ws.Range("N58").Value = -4039.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 98631.664
$ws.Range("J68").Value = 98631.664
$ws.Range("L68").Value = 98631.664
$ws.Range("N68").Value = -100129.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 98631.664
$ws.Range("J71").Value = 98631.664
$ws.Range("L71").Value = 295894.992
$ws.Range("N71").Value = -303382.992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1217.3871
$ws.Range("I105").Value = 918.0417
$ws.Range("K105").Value = 918.0417
$ws.Range("M105").Value = 828.9583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 585.1579
$ws.Range("I107").Value = 228.45454
$ws.Range("K107").Value = 228.45454
$ws.Range("M107").Value = 1691.54546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3031.475
$ws.Range("I132").Value = 2515.0386
$ws.Range("K132").Value = 7545.1158
$ws.Range("M132").Value = -5015.1158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3165.1462
$ws.Range("I134").Value = 3114.7715
$ws.Range("K134").Value = 9344.3145
$ws.Range("M134").Value = -6809.3145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2691.7896
$ws.Range("I136").Value = 1844.3
$ws.Range("J136").Value = 3633.4443
$ws.Range("K136").Value = 5532.9
$ws.Range("L136").Value = 10900.3329
$ws.Range("M136").Value = -2982.9
$ws.Range("N136").Value = -16000.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 415.5
$ws.Range("J46").Value = 549.5
$ws.Range("L46").Value = 1648.5
$ws.Range("N46").Value = -1830.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 145968.28
$ws.Range("I109").Value = 145968.28
$ws.Range("K109").Value = 437904.84
$ws.Range("M109").Value = -436864.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5661
$ws.Range("J137").Value = 1903.75
$ws.Range("L137").Value = 5711.25
$ws.Range("N137").Value = -15911.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3527.5454
$ws.Range("I141").Value = 3527.5454
$ws.Range("K141").Value = 10582.6362
$ws.Range("M141").Value = -5402.636200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 71431176
$ws.Range("I113").Value = 83335870
$ws.Range("K113").Value = 83335870
$ws.Range("M113").Value = -83333700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2587.25
$ws.Range("I132").Value = 2450
$ws.Range("K132").Value = 7350
$ws.Range("M132").Value = -4820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 506.66666
$ws.Range("I46").Value = 506.66666
$ws.Range("K46").Value = 506.66666
$ws.Range("M46").Value = -318.66666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4990.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4990.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4990.5
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -6072.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4562.08
$ws.Range("I132").Value = 4666.143
$ws.Range("K132").Value = 13998.429
$ws.Range("M132").Value = -11468.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7778.8887
$ws.Range("I136").Value = 7499.5
$ws.Range("J136").Value = 8002.4
$ws.Range("K136").Value = 22498.5
$ws.Range("L136").Value = 24007.2
$ws.Range("M136").Value = -19948.5
$ws.Range("N136").Value = -29107.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 594.625
$ws.Range("I107").Value = 643
$ws.Range("K107").Value = 1929
$ws.Range("M107").Value = -9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1982.5
$ws.Range("I126").Value = 2080.4285
$ws.Range("K126").Value = 6241.2855
$ws.Range("M126").Value = -3771.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2402.25
$ws.Range("I132").Value = 2274.88
$ws.Range("K132").Value = 6824.64
$ws.Range("M132").Value = -4294.64

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8777042
$ws.Range("I136").Value = 10421246
$ws.Range("K136").Value = 31263738
$ws.Range("M136").Value = -31261188
